$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the BTec logo inline pictures (both the "first page" and the
# "default"/primary headers) from image1.jpg -> image2.jpg.
for ($hi = 1; $hi -le 2; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists) {
        for ($si = 1; $si -le $hdr.Range.InlineShapes.Count; $si++) {
            $shp = $hdr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# Rename the Pearson Edexcel logo inline pictures (both the "first page"
# and the "default"/primary footers) from image2.png -> image1.png.
for ($fi = 1; $fi -le 2; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        for ($si = 1; $si -le $ftr.Range.InlineShapes.Count; $si++) {
            $shp = $ftr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
